# Update "想去人数" (interest count) values in both the "展览" sheet
# and the "全部类型" sheet to match the newly generated site output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5442
$ws1.Range("F4").Value  = 11761
$ws1.Range("F5").Value  = 289
$ws1.Range("F6").Value  = 597
$ws1.Range("F7").Value  = 172
$ws1.Range("F8").Value  = 287
$ws1.Range("F9").Value  = 1042
$ws1.Range("F10").Value = 99

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 5442
$ws4.Range("F7").Value  = 11761
$ws4.Range("F8").Value  = 289
$ws4.Range("F9").Value  = 597
$ws4.Range("F10").Value = 172
$ws4.Range("F13").Value = 287
$ws4.Range("F14").Value = 1042
$ws4.Range("F16").Value = 99
